# Apply the "Add files via upload" edit:
#   - Update several spec/requirement/test/result text cells on Sheet1
#     (columns E-I, rows 3-8) to the revised wording supplied by the
#     reviewer (more detail, corrected grammar, added "Not Met" result).
#   - Adjust the row heights for rows 4, 7 and 8 to fit the longer,
#     wrapped text.
#   - Move the active selection to I7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 : "Can estimate ..." requirement gains "the number of" ---
$ws.Range("E3").Value = "Can estimate the number of fishes fishes in the area."

# --- Row 4 : fish-food dispenser spec / test wording expanded ---
$ws.Range("F4").Value = "Dispenses I kg of fish food 1 time when button is pressed"
$ws.Range("H4").Value = "Bottle rotates and dispenses the fish food ( the system should be modify to be able to dispense more food mutiple times)"

# --- Row 6 : maintenance spec replaced with wiring/battery note ---
$ws.Range("F6").Value = "The wiring should be simple and clear with labers and the recharging of batt"

# --- Row 7 : floatation spec + waver/shake test reworded, extra detail ---
$ws.Range("F7").Value = "Floats stabaly on water surface when stationary and moving, the main body should be 10mm above water level with addition load of 2kg."
$ws.Range("H7").Value = "Does not waver or shake when moving or stationary with a low centre of gravity, the main body is above the water level."

# --- Row 8 : IR control spec/test reworded and result now "Not Met" ---
$ws.Range("F8").Value = "Using IR sensor and remote to control craft's movement and direction ( control disstance longer than 100m)"
$ws.Range("H8").Value = "Pressing buttons will make the motors rotate in the intended directions. However, the control distance is just around 5m. Should use ESP-32)"
$ws.Range("I8").Value = "Not Met"

# --- Row heights grown to fit the newly wrapped, longer text ---
$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 57.6
$ws.Rows.Item(8).RowHeight = 43.2

# --- Active cell / selection moved to I7 ---
$ws.Range("I7").Select()
